$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '48.186.17'
$ws.Cells.Item(2, 5).Value = '  +2.06%  '

$ws.Cells.Item(3, 4).Value = '2.531.66'
$ws.Cells.Item(3, 5).Value = '  +1.45%  '

$ws.Cells.Item(4, 5).Value = '  +0.01%  '

$ws.Cells.Item(5, 4).Value = '''324.39'
$ws.Cells.Item(5, 5).Value = '  +0.35%  '

$ws.Cells.Item(6, 4).Value = '''108.99'
$ws.Cells.Item(6, 5).Value = '  -0.06%  '

$ws.Cells.Item(7, 5).Value = '  +0.49%  '

$ws.Cells.Item(8, 5).Value = '  +0.04%  '

$ws.Cells.Item(9, 4).Value = '''0.558'
$ws.Cells.Item(9, 5).Value = '  +4.09%  '

$ws.Cells.Item(10, 4).Value = '''40.65'
$ws.Cells.Item(10, 5).Value = '  +3.95%  '

$ws.Cells.Item(11, 4).Value = '''20.52'
$ws.Cells.Item(11, 5).Value = '  +11.31%  '

$ws.Cells.Item(12, 4).Value = '''0.0828'
$ws.Cells.Item(12, 5).Value = '  +1.75%  '

$ws.Cells.Item(14, 5).Value = '  +1.47%  '

$ws.Cells.Item(15, 4).Value = '2.928.77'

$ws.Cells.Item(16, 4).Value = '2.537.84'
$ws.Cells.Item(16, 5).Value = '  +1.34%  '

$ws.Cells.Item(17, 5).Value = '  +1.02%  '

$ws.Cells.Item(18, 4).Value = '48.045.70'
$ws.Cells.Item(18, 5).Value = '  +1.89%  '

$ws.Cells.Item(19, 4).Value = '''13.30'
$ws.Cells.Item(19, 5).Value = '  +4.08%  '

$ws.Cells.Item(20, 4).Value = '''6.66'
$ws.Cells.Item(20, 5).Value = '  +0.33%  '

$ws.Cells.Item(21, 5).Value = '  +1.18%  '

$ws.Cells.Item(22, 5).Value = '  -1.08%  '

$ws.Cells.Item(23, 4).Value = '''72.31'
$ws.Cells.Item(23, 5).Value = '  +2.24%  '

$ws.Cells.Item(24, 4).Value = '''269.72'
$ws.Cells.Item(24, 5).Value = '  +8.93%  '

$ws.Cells.Item(25, 5).Value = '  -0.24%  '

$ws.Cells.Item(26, 4).Value = '''26.25'

$ws.Cells.Item(28, 4).Value = '''10.17'
$ws.Cells.Item(28, 5).Value = '  +0.85%  '

$ws.Cells.Item(29, 5).Value = '  +2.44%  '

$ws.Cells.Item(30, 5).Value = '  -3.72%  '

$ws.Cells.Item(31, 4).Value = '''35.72'
$ws.Cells.Item(31, 5).Value = '  +1.08%  '

$ws.Cells.Item(32, 5).Value = '  -0.43%  '

$ws.Cells.Item(33, 4).Value = '''19.89'
$ws.Cells.Item(33, 5).Value = '  -0.65%  '

$ws.Cells.Item(34, 4).Value = '''5.43'
$ws.Cells.Item(34, 5).Value = '  +0.13%  '

$ws.Cells.Item(35, 5).Value = '  +0.01%  '

$ws.Cells.Item(36, 5).Value = '  +1.05%  '

$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(37, 4).Value = '''4.78'
$ws.Cells.Item(37, 5).Value = '  +1.78%  '

$ws.Cells.Item(38, 2).Value = 'ARBITRUM'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(38, 4).Value = '''2.00'
$ws.Cells.Item(38, 5).Value = '  +0.89%  '

$ws.Cells.Item(39, 5).Value = '  +1.01%  '

$ws.Cells.Item(40, 5).Value = '  +0.09%  '

$ws.Cells.Item(41, 4).Value = '''22.43'
$ws.Cells.Item(41, 5).Value = '  +6.21%  '

$ws.Cells.Item(42, 4).Value = '''119.55'
$ws.Cells.Item(42, 5).Value = '  -1.97%  '

$ws.Cells.Item(43, 5).Value = '  -1.24%  '

$ws.Cells.Item(44, 5).Value = '  +0.80%  '

$ws.Cells.Item(45, 4).Value = '2.013.90'
$ws.Cells.Item(45, 5).Value = '  +1.01%  '

$ws.Cells.Item(46, 4).Value = '''3.14'
$ws.Cells.Item(46, 5).Value = '  +3.07%  '

$ws.Cells.Item(47, 5).Value = '  +0.01%  '

$ws.Cells.Item(48, 5).Value = '  +4.98%  '

$ws.Cells.Item(49, 4).Value = '''9.17'
$ws.Cells.Item(49, 5).Value = '  +0.45%  '

$ws.Cells.Item(50, 5).Value = '  +1.38%  '

$ws.Cells.Item(51, 4).Value = '''79.80'
$ws.Cells.Item(51, 5).Value = '  +2.30%  '

